# "+ bổ sung lỗi" -- update the error list on the "Quản lý giao dịch _ Nạp tiền"
# item (row 10-12): the two separate notes in D10/D11 are consolidated into a
# single, updated note in D10, and the now-unused D11/D12 sub-bullet cells are
# cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the note text in D10 with the new bug description.
$ws.Range("D10").Value = "_ nhập số thẻ đúng (copy từ csdl tình trạng =1, =0 luôn) nhưng cứ báo tài khoản không đúng."

# The note now wraps onto two lines, so the row needs to be taller.
$ws.Rows.Item(10).RowHeight = 33

# D11 previously held a second note; clear its contents but keep the cell's
# existing border/format.
$ws.Range("D11").ClearContents()

# D12 was already empty and had no note of its own - fully clear it (contents
# and formatting) so it drops back to the sheet's default, unused cell.
$ws.Range("D12").Clear()

# Move the active selection to D12, matching where the edit left off.
[void]$ws.Range("D12").Select()
